# "rimosso header dalla prima pagina" - apply the 1 MASTER.docx edit:
#   1) Increase the last section's page margins (the "Media List" section),
#      which now has extra top space since its header text is gone from page 1.
#   2) Refresh the shared header's two VML pict shapes (new anchorIds) and
#      update the cached STYLEREF field result from "Media List" to
#      "Vision Statement".

$d = $word.ActiveDocument

# --- 1) Page margins (last / body-level section) ---------------------------
# 567 twips -> 1417/1134 twips == 28.35pt -> 70.85pt / 56.7pt (1 pt = 20 twips)
$lastSection = $d.Sections.Last
$lastSection.PageSetup.TopMargin    = 70.85
$lastSection.PageSetup.RightMargin  = 56.7
$lastSection.PageSetup.BottomMargin = 56.7
$lastSection.PageSetup.LeftMargin   = 56.7

# --- 2) Shared header (word/header1.xml) ------------------------------------
# The legacy VML text-box shapes inside the header aren't addressable through
# Shapes()/TextFrame on this host, so rebuild the header paragraph verbatim
# (same structure Word itself emits) with only the three changed tokens:
#   w14:anchorId 19A5B48A -> 12116593
#   w14:anchorId 050BC11E -> 40728461
#   field result "Media List" -> "Vision Statement"
$headerXml = '<w:hdr xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid wp14"><w:p w14:paraId="7B766B8F" w14:textId="77777777" w:rsidR="001B03CE" w:rsidRDefault="002C2B6E"><w:pPr><w:pStyle w:val="Header"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="12116593"><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="Text Box 218" o:spid="_x0000_s2050" type="#_x0000_t202" style="position:absolute;left:0;text-align:left;margin-left:56.7pt;margin-top:28.5pt;width:481.9pt;height:13.8pt;z-index:2;visibility:visible;mso-wrap-style:square;mso-width-percent:1000;mso-height-percent:0;mso-left-percent:-10001;mso-top-percent:-10001;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:page;mso-width-percent:1000;mso-height-percent:0;mso-left-percent:-10001;mso-top-percent:-10001;mso-width-relative:margin;mso-height-relative:page;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQCaUtwRswIAALcFAAAOAAAAZHJzL2Uyb0RvYy54bWysVNtu2zAMfR+wfxD07vpSJY2NOkUbx8OA&#xA;7gK0+wBFlmNhtuRJSpxu2L+PkpM0aTFg2OYHQ6KoQx7yiNc3u65FW66NUDLH8UWEEZdMVUKuc/zl&#xA;sQxmGBlLZUVbJXmOn7jBN/O3b66HPuOJalRbcY0ARJps6HPcWNtnYWhYwztqLlTPJRzWSnfUwlav&#xA;w0rTAdC7NkyiaBoOSle9VowbA9ZiPMRzj1/XnNlPdW24RW2OITfr/9r/V+4fzq9ptta0bwTbp0H/&#xA;IouOCglBj1AFtRRttHgF1QmmlVG1vWCqC1VdC8Y9B2ATRy/YPDS0554LFMf0xzKZ/wfLPm4/aySq&#xA;HCcxtErSDpr0yHcW3akdcjao0NCbDBwfenC1OziATnu2pr9X7KtBUi0aKtf8Vms1NJxWkGHsboYn&#xA;V0cc40BWwwdVQSC6scoD7WrdufJBQRCgQ6eejt1xyTAwTlJyOY3giMFZfBXN4okPQbPD7V4b+46r&#xA;DrlFjjV036PT7b2xLhuaHVxcMKlK0bZeAa08M4DjaIHYcNWduSx8Q3+kUbqcLWckIMl0GZCoKILb&#xA;ckGCaRlfTYrLYrEo4p8ubkyyRlQVly7MQVwx+bPm7WU+yuIoL6NaUTk4l5LR69Wi1WhLQdyl//YF&#xA;OXELz9PwRQAuLyjFCYnukjQop7OrgJRkEqRQ4CCK07t0GpGUFOU5pXsh+b9TQkOO00kyGcX0W26R&#xA;/15zo1knLIyPVnQ5nh2daOYkuJSVb62loh3XJ6Vw6T+XAtp9aLQXrNPoqFa7W+0Axal4paonkK5W&#xA;oCwQIcw8WDRKf8dogPmRY/NtQzXHqH0vQf5pTIgbOH4DC31qXR2sVDKAyDGzGqNxs7DjeNr0Wqwb&#xA;iHF4arfwWErhdfycz/6JwXTwdPaTzI2f0733ep63818AAAD//wMAUEsDBBQABgAIAAAAIQBczPU/&#xA;2wAAAAQBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI9BS8NAEIXvgv9hGcGb3VhLMDGbIoIepCpGaa/T&#xA;7JgEs7Mxu23Tf+/oRS8Djze8971iOble7WkMnWcDl7MEFHHtbceNgfe3+4trUCEiW+w9k4EjBViW&#xA;pycF5tYf+JX2VWyUhHDI0UAb45BrHeqWHIaZH4jF+/CjwyhybLQd8SDhrtfzJEm1w46locWB7lqq&#xA;P6udk5L1Ex6fk5V7qR+/sofNqqkWi8aY87Pp9gZUpCn+PcMPvqBDKUxbv2MbVG9AhsTfK152lYrc&#xA;GpinGeiy0P/hy28AAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAA&#xA;AAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAA&#xA;AAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAmlLcEbMCAAC3BQAADgAA&#xA;AAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAXMz1P9sAAAAEAQAA&#xA;DwAAAAAAAAAAAAAAAAANBQAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAABUGAAAAAA==&#xA;" o:allowincell="f" filled="f" stroked="f"><v:textbox style="mso-next-textbox:#Text Box 218;mso-fit-shape-to-text:t" inset=",0,,0"><w:txbxContent><w:p w14:paraId="2270C80B" w14:textId="77777777" w:rsidR="001B03CE" w:rsidRDefault="002C2B6E" w:rsidP="004C1716"><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> STYLEREF  "Heading 1"  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00926B59"><w:rPr><w:noProof/></w:rPr><w:t>Vision Statement</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:txbxContent></v:textbox><w10:wrap anchorx="margin" anchory="margin"/></v:shape></w:pict></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="40728461"><v:shape id="Text Box 219" o:spid="_x0000_s2049" type="#_x0000_t202" style="position:absolute;left:0;text-align:left;margin-left:0;margin-top:28.5pt;width:56.7pt;height:13.8pt;z-index:1;visibility:visible;mso-wrap-style:square;mso-width-percent:1000;mso-height-percent:0;mso-left-percent:-10001;mso-top-percent:-10001;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:page;mso-width-percent:1000;mso-height-percent:0;mso-left-percent:-10001;mso-top-percent:-10001;mso-width-relative:left-margin-area;mso-height-relative:page;v-text-anchor:middle" o:allowincell="f" fillcolor="#f79e10" stroked="f"><v:textbox style="mso-fit-shape-to-text:t" inset=",0,,0"><w:txbxContent><w:p w14:paraId="1835FFD2" w14:textId="77777777" w:rsidR="001B03CE" w:rsidRPr="001B03CE" w:rsidRDefault="001B03CE"><w:pPr><w:jc w:val="right"/><w:rPr><w:color w:val="FFFFFF"/></w:rPr></w:pPr><w:r w:rsidRPr="001B03CE"><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r><w:r w:rsidRPr="001B03CE"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidRPr="001B03CE"><w:rPr><w:noProof/><w:color w:val="FFFFFF"/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidRPr="001B03CE"><w:rPr><w:noProof/><w:color w:val="FFFFFF"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:txbxContent></v:textbox><w10:wrap anchorx="page" anchory="margin"/></v:shape></w:pict></w:r></w:p></w:hdr>'

$header = $d.Sections.Item(1).Headers.Item(1)
$header.Range.InsertXML($headerXml) | Out-Null
